$wb = $excel.ActiveWorkbook

# --- AddOffercode sheet: update offer code T146 -> T147 ---
$wsAdd = $wb.Worksheets.Item("AddOffercode")
$wsAdd.Range("A2").Value = "T147"
$wsAdd.Range("B2").Value = 20
$wsAdd.Range("D2").Value = "Offer code T147 description"

# --- EditOfferCode sheet: move selection from B8 to A2 ---
$wsEdit = $wb.Worksheets.Item("EditOfferCode")
$wsEdit.Activate()
$wsEdit.Range("A2").Select()

# --- AddOffercode becomes the active/selected tab, with B2 selected ---
$wsAdd.Activate()
$wsAdd.Range("B2").Select()
